$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Bmp6-Bmpr1a LR-pair rows: ECs is now also included as a
# sending cluster, giving a full 3x3 (ECs/FAPs/sCs) x (ECs/FAPs/sCs) matrix.
$data = @(
    @("ECs", "Bmp6", "Bmpr1a", "ECs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 9.918858999999999, 29.756577, 0.160764128269069, 0.160764128269069, 315.672657023359, 2841.053913210231, 0.09748107950400781, 0.09748107950400781),
    @("ECs", "Bmp6", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 29.20351433333333, 87.61054300000001, 0.47332838627826, 0.4733283862782601, 929.4164746190145, 8364.74827157113, 0.2870078204079822, 0.2870078204079823),
    @("ECs", "Bmp6", "Bmpr1a", "sCs", 2, 0.6666666666666666, 31.825501, 95.47650300000001, 0.6063608875535647, 0.6063608875535647, 3, 1, 22.575837, 67.72751099999999, 0.3659074854526709, 0.3659074854526709, 718.4873230193369, 6466.385907174033, 0.2218719876415746, 0.2218719876415746),
    @("FAPs", "Bmp6", "Bmpr1a", "ECs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 9.918858999999999, 29.756577, 0.160764128269069, 0.160764128269069, 48.06957460709265, 432.6261714638339, 0.0148440921940575, 0.0148440921940575),
    @("FAPs", "Bmp6", "Bmpr1a", "FAPs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 29.20351433333333, 87.61054300000001, 0.47332838627826, 0.4733283862782601, 141.5284269123562, 1273.755842211206, 0.04370458932367924, 0.04370458932367924),
    @("FAPs", "Bmp6", "Bmpr1a", "sCs", 3, 1, 4.846280666666666, 14.538842, 0.09233460445363234, 0.09233460445363234, 3, 1, 22.575837, 67.72751099999999, 0.3659074854526709, 0.3659074854526709, 109.408842386918, 984.6795814822618, 0.03378592293589559, 0.0337859229358956),
    @("sCs", "Bmp6", "Bmpr1a", "ECs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 9.918858999999999, 29.756577, 0.160764128269069, 0.160764128269069, 156.85971269511, 1411.73741425599, 0.04843895657100371, 0.04843895657100372),
    @("sCs", "Bmp6", "Bmpr1a", "FAPs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 29.20351433333333, 87.61054300000001, 0.47332838627826, 0.4733283862782601, 461.83284468649, 4156.49560217841, 0.1426159765465985, 0.1426159765465986),
    @("sCs", "Bmp6", "Bmpr1a", "sCs", 3, 1, 15.81429, 47.44287, 0.301304507992803, 0.3013045079928031, 3, 1, 22.575837, 67.72751099999999, 0.3659074854526709, 0.3659074854526709, 357.0208333107299, 3213.18749979657, 0.1102495748752007, 0.1102495748752008),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}